# Auto-generated script: appends the new "Knarot" section + reference list after
# the final paragraph ("BILAGA 1 - Fridlysta arter"). Plain text for every new
# paragraph is written first (so no paragraph can inherit italic formatting from
# a previous one); italics are then applied in a second pass over recorded
# absolute ranges. Finally, the header date is updated.

$d = $word.ActiveDocument

$anchor = $d.Paragraphs($d.Paragraphs.Count)
$italicStarts = New-Object System.Collections.ArrayList
$italicEnds = New-Object System.Collections.ArrayList

# --- New paragraph 1 of 13 (style: Heading1) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Heading1"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'

# --- New paragraph 2 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'

# --- New paragraph 3 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
[void]$italicStarts.Add($pStart + 34)
[void]$italicEnds.Add($pStart + 116)
[void]$italicStarts.Add($pStart + 278)
[void]$italicEnds.Add($pStart + 483)
[void]$italicStarts.Add($pStart + 490)
[void]$italicEnds.Add($pStart + 608)

# --- New paragraph 4 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
[void]$italicStarts.Add($pStart + 205)
[void]$italicEnds.Add($pStart + 1070)

# --- New paragraph 5 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'

# --- New paragraph 6 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'

# --- New paragraph 7 of 13 (style: Heading2) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Heading2"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Referenser - knärot'

# --- New paragraph 8 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
[void]$italicStarts.Add($pStart + 33)
[void]$italicEnds.Add($pStart + 113)

# --- New paragraph 9 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
[void]$italicStarts.Add($pStart + 62)
[void]$italicEnds.Add($pStart + 176)

# --- New paragraph 10 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
[void]$italicStarts.Add($pStart + 117)
[void]$italicEnds.Add($pStart + 207)

# --- New paragraph 11 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
[void]$italicStarts.Add($pStart + 54)
[void]$italicEnds.Add($pStart + 121)

# --- New paragraph 12 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
[void]$italicStarts.Add($pStart + 22)
[void]$italicEnds.Add($pStart + 57)

# --- New paragraph 13 of 13 (style: Normal) ---
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "Normal"
$pStart = $anchor.Range.Start
$anchor.Range.Text = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
[void]$italicStarts.Add($pStart + 25)
[void]$italicEnds.Add($pStart + 61)

# --- Second pass: apply italics now that every paragraph mark already exists ---
for ($i = 0; $i -lt $italicStarts.Count; $i++) {
    $ir = $d.Range($italicStarts[$i], $italicEnds[$i])
    $ir.Font.Italic = $true
}

# --- Update the date in the first-page header (2023-09-13 -> 2023-09-15) ---
foreach ($sec in $d.Sections) {
    foreach ($hfType in @(1, 2, 3)) {
        $hdr = $sec.Headers.Item($hfType)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false,
                                     $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}

Write-Output "done"